$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OES Sheet")
$ws.Rows("1:5").Delete()
